# Update "想去人数" (F column) figures for several rows across sheets,
# and for two rows that became sellable, replace the "不可售" text in the
# "最低票价" (G column) with a numeric price of 39.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 872
$ws1.Range("F3").Value  = 997
$ws1.Range("F4").Value  = 781
$ws1.Range("F6").Value  = 440
$ws1.Range("F7").Value  = 679
$ws1.Range("F8").Value  = 154
$ws1.Range("F9").Value  = 1280
$ws1.Range("F10").Value = 707
$ws1.Range("F12").Value = 544
$ws1.Range("F13").Value = 182
$ws1.Range("F14").Value = 35
$ws1.Range("F15").Value = 939
$ws1.Range("F16").Value = 10
$ws1.Range("G16").Value = 39
$ws1.Range("F17").Value = 398
$ws1.Range("F18").Value = 373
$ws1.Range("F19").Value = 91
$ws1.Range("F21").Value = 139
$ws1.Range("F22").Value = 629
$ws1.Range("F24").Value = 970
$ws1.Range("F25").Value = 13

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 27
$ws2.Range("F6").Value = 189
$ws2.Range("F7").Value = 238
$ws2.Range("F9").Value = 29

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 872
$ws4.Range("F5").Value  = 997
$ws4.Range("F6").Value  = 781
$ws4.Range("F8").Value  = 440
$ws4.Range("F9").Value  = 679
$ws4.Range("F10").Value = 154
$ws4.Range("F11").Value = 1280
$ws4.Range("F12").Value = 707
$ws4.Range("F14").Value = 27
$ws4.Range("F16").Value = 544
$ws4.Range("F18").Value = 182
$ws4.Range("F19").Value = 35
$ws4.Range("F20").Value = 940
$ws4.Range("F21").Value = 189
$ws4.Range("F22").Value = 10
$ws4.Range("G22").Value = 39
$ws4.Range("F23").Value = 398
$ws4.Range("F24").Value = 373
$ws4.Range("F25").Value = 91
$ws4.Range("F26").Value = 238
$ws4.Range("F29").Value = 29
$ws4.Range("F33").Value = 139
$ws4.Range("F34").Value = 629
$ws4.Range("F36").Value = 970
$ws4.Range("F37").Value = 13
